# Refresh the cryptocurrency price/volume snapshot (symbol list update run by
# the scheduled GitHub Actions job). Only the "Price" (D) and "Volume(1h)"
# (E) columns move; every other column is left untouched.
#
# The source cells store these figures as literal text (e.g. "305.90",
# "0.86%") rather than numeric values, so each target cell is first switched
# to the Text number format before the new literal is written. That mirrors
# how the figures ended up as text in the workbook in the first place and
# keeps Excel from reinterpreting "306.03" as a number or "0.98%" as a
# percentage.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2";  Value = "306.03" }
    @{ Cell = "E2";  Value = "0.98%" }
    @{ Cell = "E3";  Value = "-1.78%" }
    @{ Cell = "D4";  Value = "5.045" }
    @{ Cell = "E4";  Value = "0.56%" }
    @{ Cell = "D5";  Value = "0.07897" }
    @{ Cell = "E5";  Value = "2.00%" }
    @{ Cell = "D6";  Value = "2.267" }
    @{ Cell = "E6";  Value = "8.05%" }
    @{ Cell = "D7";  Value = "7.998" }
    @{ Cell = "E7";  Value = "0.04%" }
    @{ Cell = "D8";  Value = "4.147" }
    @{ Cell = "E8";  Value = "2.31%" }
    @{ Cell = "D9";  Value = "0.9279" }
    @{ Cell = "D10"; Value = "0.09803" }
    @{ Cell = "E10"; Value = "-0.13%" }
    @{ Cell = "D11"; Value = "0.1868" }
    @{ Cell = "E11"; Value = "0.40%" }
    @{ Cell = "D12"; Value = "0.09004" }
    @{ Cell = "E12"; Value = "5.08%" }
    @{ Cell = "D13"; Value = "0.03744" }
    @{ Cell = "E13"; Value = "3.80%" }
    @{ Cell = "D14"; Value = "0.09912" }
    @{ Cell = "E14"; Value = "-0.64%" }
    @{ Cell = "D15"; Value = "0.001441" }
    @{ Cell = "E15"; Value = "-2.32%" }
    @{ Cell = "D16"; Value = "0.005716" }
    @{ Cell = "E16"; Value = "-0.94%" }
    @{ Cell = "D17"; Value = "3.462" }
    @{ Cell = "E17"; Value = "-0.21%" }
    @{ Cell = "E18"; Value = "4.19%" }
    @{ Cell = "E19"; Value = "-1.98%" }
    @{ Cell = "E20"; Value = "-1.11%" }
    @{ Cell = "D21"; Value = "5.081" }
    @{ Cell = "E21"; Value = "2.18%" }
    @{ Cell = "D22"; Value = "0.2250" }
    @{ Cell = "E22"; Value = "1.57%" }
    @{ Cell = "D23"; Value = "0.04577" }
    @{ Cell = "E23"; Value = "-0.56%" }
    @{ Cell = "D24"; Value = "0.001233" }
    @{ Cell = "E24"; Value = "-0.54%" }
    @{ Cell = "D25"; Value = "0.004778" }
    @{ Cell = "E25"; Value = "-6.56%" }
    @{ Cell = "D26"; Value = "0.0001301" }
    @{ Cell = "E26"; Value = "-7.54%" }
    @{ Cell = "D39"; Value = "0.01932" }
    @{ Cell = "E39"; Value = "9.04%" }
    @{ Cell = "D40"; Value = "0.04944" }
    @{ Cell = "E40"; Value = "6.12%" }
    @{ Cell = "E41"; Value = "1.39%" }
    @{ Cell = "E42"; Value = "0.01%" }
    @{ Cell = "D43"; Value = "0.007811" }
    @{ Cell = "E43"; Value = "2.23%" }
    @{ Cell = "D44"; Value = "0.002192" }
    @{ Cell = "E44"; Value = "1.23%" }
    @{ Cell = "E45"; Value = "15.28%" }
    @{ Cell = "D46"; Value = "0.00006143" }
    @{ Cell = "E46"; Value = "-2.52%" }
    @{ Cell = "E47"; Value = "-0.57%" }
    @{ Cell = "D48"; Value = "51.76" }
    @{ Cell = "E48"; Value = "41.93%" }
    @{ Cell = "D49"; Value = "0.001801" }
    @{ Cell = "E49"; Value = "-10.51%" }
    @{ Cell = "D50"; Value = "0.00002102" }
    @{ Cell = "E50"; Value = "-0.57%" }
    @{ Cell = "D51"; Value = "0.0002002" }
    @{ Cell = "E51"; Value = "-0.57%" }
)

foreach ($update in $updates) {
    $cell = $ws.Range($update.Cell)
    $cell.NumberFormat = "@"
    $cell.Value = $update.Value
}
